# "Added React Native project"
#
# 1) Correct a typo in the existing Xamarin row: "POC" -> "Porototipo".
# 2) Append a new tracked-hours row for the React Native prototype.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# 1. Fix existing text in C5 ("POC applicativo..." -> "Porototipo applicativo...")
$ws.Range("C5").Value = "Porototipo applicativo mobile con Xamarin"

# 2. Add the new row (row 6) just below the last used row.
#    Copy the formatting from the row above so date/number/text styles match.
$ws.Range("A5:C5").Copy() | Out-Null
$ws.Range("A6:C6").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$wb.Application.CutCopyMode = $false

$ws.Range("A6").Value = "7/29/2022"
$ws.Range("B6").Value = 3
$ws.Range("C6").Value = "Porototipo applicativo mobile con React Native"

# Move the selection down to the next empty row and scroll the sheet so
# column C (Descrizione) stays in view, matching the saved workbook view.
$excel.ActiveWindow.ScrollColumn = 3
$ws.Range("C7").Select() | Out-Null
